$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bwi")

# Update the destination value in A2 from "New York, NY, USA" to "Chicago, IL, USA"
$ws.Range("A2").Value = "Chicago, IL, USA"

# Update the selected cell to B4
$ws.Range("B4").Select()
